# DBMapperTests.xlsx edit script
# Simulates a refreshed DB query pull (fewer records), a TestId -> TestID
# column rename, and the associated bookkeeping (defined names, table,
# comments, selections) that follow from it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Example2": DBListFetch-driven dump (with extra G:Q helper cells)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Example2")
$ws2.UsedRange.ClearContents()

# Row 1: new formula row (query status / record count)
$ws2.Range("A1").Formula = '=_xll.DBListFetch("select * from ORE..Testtable","",A2,,,TRUE,TRUE,TRUE)'

# Row 2: header row (was row 1)
$ws2.Range("A2").Value2 = "TestID"
$ws2.Range("B2").Value2 = "TestStr"
$ws2.Range("C2").Value2 = "TestDate"
$ws2.Range("D2").Value2 = "TestNum"
$ws2.Range("E2").Value2 = "ignoredColumn"

# Data rows 3-10 (8 records)
$ws2.Range("A3").Value2 = 1
$ws2.Range("B3").Value2 = "testtest"
$ws2.Range("C3").Value2 = 32874
$ws2.Range("D3").Value2 = 123.9
$ws2.Range("E3").Value2 = "cvb"

$ws2.Range("A4").Value2 = 2
$ws2.Range("B4").Value2 = "testtesttest"
$ws2.Range("C4").Value2 = 43586
$ws2.Range("D4").Value2 = 147
$ws2.Range("E4").Value2 = "cvb"

$ws2.Range("A5").Value2 = 3
$ws2.Range("B5").Value2 = "sdfsdf"
$ws2.Range("C5").Value2 = 43777
$ws2.Range("D5").Value2 = 456.25
$ws2.Range("E5").Value2 = "ghfgh"

$ws2.Range("A6").Value2 = 4
$ws2.Range("B6").Value2 = "rewrwer"
$ws2.Range("C6").Value2 = 43685
$ws2.Range("D6").Value2 = 478.32

$ws2.Range("A7").Value2 = 5
$ws2.Range("B7").Value2 = "werwer"
$ws2.Range("C7").Value2 = 43717
$ws2.Range("D7").Value2 = 654

$ws2.Range("A8").Value2 = 6
$ws2.Range("B8").Value2 = "zrtzrtz"
$ws2.Range("C8").Value2 = 43685
$ws2.Range("D8").Value2 = 457.5

$ws2.Range("A9").Value2 = 7
$ws2.Range("B9").Value2 = "rtzrtzrtzrtz"
$ws2.Range("C9").Value2 = 43717
$ws2.Range("D9").Value2 = 5

$ws2.Range("A10").Value2 = 8
$ws2.Range("B10").Value2 = "ertert"
$ws2.Range("C10").Value2 = 43685

# Rows 11-12 stay empty (formatting leftovers only in the original file)

# Helper block G3:Q3 (was G2:Q2) - blank formatted cells
$ws2.Range("G3:Q3").ClearContents()

# Column A width nudge (matches bestFit after header text change)
$ws2.Columns.Item(1).ColumnWidth = 6.5703125

# ---------------------------------------------------------------------
# Sheet "Example3": DBSetQuery-driven dump, backs the TestTable table
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Example3")
$ws3.UsedRange.ClearContents()

$ws3.Range("A1").Formula = '=_xll.DBSetQuery("Select * FROM ORE..TestTable","",A2)'

$ws3.Range("A2").Value2 = "TestID"
$ws3.Range("B2").Value2 = "TestStr"
$ws3.Range("C2").Value2 = "TestDate"
$ws3.Range("D2").Value2 = "TestNum"

$ws3.Range("A3").Value2 = 1
$ws3.Range("B3").Value2 = "testtest"
$ws3.Range("C3").Value2 = 32874
$ws3.Range("D3").Value2 = 123.9

$ws3.Range("A4").Value2 = 2
$ws3.Range("B4").Value2 = "testtesttest"
$ws3.Range("C4").Value2 = 43586
$ws3.Range("D4").Value2 = 147

$ws3.Range("A5").Value2 = 3
$ws3.Range("B5").Value2 = "sdfsdf"
$ws3.Range("C5").Value2 = 43777
$ws3.Range("D5").Value2 = 456.25

$ws3.Range("A6").Value2 = 4
$ws3.Range("B6").Value2 = "rewrwer"
$ws3.Range("C6").Value2 = 43685
$ws3.Range("D6").Value2 = 478.32

$ws3.Range("A7").Value2 = 5
$ws3.Range("B7").Value2 = "werwer"
$ws3.Range("C7").Value2 = 43717
$ws3.Range("D7").Value2 = 654

$ws3.Range("A8").Value2 = 6
$ws3.Range("B8").Value2 = "zrtzrtz"
$ws3.Range("C8").Value2 = 43685
$ws3.Range("D8").Value2 = 457.5

$ws3.Range("A9").Value2 = 7
$ws3.Range("B9").Value2 = "rtzrtzrtzrtz"
$ws3.Range("C9").Value2 = 43717
$ws3.Range("D9").Value2 = 5

$ws3.Range("A10").Value2 = 8
$ws3.Range("B10").Value2 = "ertert"
$ws3.Range("C10").Value2 = 43685

$ws3.Columns.Item(1).ColumnWidth = 8.85546875

# Resize+rename the backing table (TestId -> TestID, A1:D11 -> A2:D10)
$lo = $ws3.ListObjects.Item(1)
$lo.Resize($ws3.Range("A2:D10"))
$lo.ListColumns.Item(1).Name = "TestID"

# ---------------------------------------------------------------------
# Workbook-level defined names
# ---------------------------------------------------------------------
$names = $wb.Names

$n = $names.Item("DBFsource43642653946169")
$n.RefersTo = "=Example2!`$G`$2"

$names.Add("DBFsource436715641890509", "=Example2!`$A`$1")
$nHidden = $names.Item("DBFsource436715641890509")
$nHidden.Visible = $false

$n2 = $names.Item("DBFtarget43642653946169")
$n2.Name = "DBFtarget436715641890509"
$n2.RefersTo = "=Example2!`$A`$2:`$D`$10"

$n3 = $names.Item("DBMapper")
$n3.RefersTo = "=Example2!`$A`$2"

$n4 = $names.Item("DBMapperDataRange")
$n4.RefersTo = "=TestTable[[#Headers],[TestID]]"

$n5 = $names.Item("OEBFADBTVI00_ORE_TestTable")
$n5.RefersTo = "=Example3!`$A`$2:`$D`$10"

# ---------------------------------------------------------------------
# Comments: shift from row1 to row2 on Example2 / Example3
# ---------------------------------------------------------------------
$c2 = $ws2.Comments.Item(1)
$c2Text = $c2.Text()
$c2.Delete()
$ws2.Range("A2").AddComment($c2Text)

$c3 = $ws3.Comments.Item(1)
$c3Text = $c3.Text()
$c3.Delete()
$ws3.Range("A2").AddComment($c3Text)

# ---------------------------------------------------------------------
# Selections / active sheet (Example3 ends up active, tab index 2)
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B7").Select()

$ws3.Activate()
$ws3.Range("D10").Select()
